$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B=29.12242425684365; C=23.02353161545612; D=5.457562313420514; E=29.30942974478529; F=42.83252900737337; G=2.069037592771776; H=3.219152413621614; I=3.413434991520969; P=13.37876734696289 }
    3 = @{ B=27.1374243989585; C=21.40733941091886; D=5.337267986424786; E=27.25656006138376; F=40.45473511116975; G=2.078198967281577; H=2.859703050369023; I=3.074111969351026; P=13.44169929310584 }
    4 = @{ B=25.87105930567509; C=20.36552448339359; D=5.260392023805906; E=25.9310335883985; F=38.9386664220199; G=2.083948925629955; H=2.633168838164859; I=2.861811536839314; P=13.48005957020183 }
    5 = @{ B=25.32578769390431; C=19.93426239784606; D=5.224182377235195; E=25.37362898793707; F=38.28229444321997; G=2.086343736601319; H=2.538891211510318; I=2.774232242323455; P=13.49076358101514 }
    6 = @{ B=25.21926499381419; C=19.86909335459384; D=5.213083032900629; E=25.27925579919151; F=38.14214131045867; G=2.086765773327855; H=2.522620236960007; I=2.759813657969449; P=13.48663518098441 }
    7 = @{ B=25.82386594742595; C=20.37891803677698; D=5.24621841167165; E=25.92137325416789; F=38.84853422600711; G=2.084041586545941; H=2.630487630487425; I=2.860669597260886; P=13.46411645950251 }
    8 = @{ B=28.39101441007527; C=22.49958675012746; D=5.398999121240402; E=28.61202614888293; F=41.92355994330617; G=2.072250082095517; H=3.094540212180184; I=3.296299062931547; P=13.3793457196131 }
    9 = @{ B=33.0895634734422; C=26.24400022884295; D=5.698493296056921; E=33.386799994948; F=47.63182860177298; G=2.049914329690132; H=3.964171270026597; I=4.12457982399691; P=13.24067711716932 }
    10 = @{ B=36.0434723756846; C=28.63441750164197; D=5.830995088903073; E=35.71059044858454; F=51.13791178948935; G=2.034521241778429; H=4.534438610626464; I=4.698664340282644; P=13.06500675699586 }
    11 = @{ B=36.02269999258484; C=28.45708531374927; D=5.314150617902084; E=29.16197763703897; F=49.38345475826942; G=2.03254108020656; H=4.85817995592064; I=4.774015261544365; P=12.41367954742661 }
    12 = @{ B=35.45132614581257; C=27.80348711356223; D=4.902634528648329; E=23.22781200176598; F=47.2578856750118; G=2.033596121678061; H=5.582613299285839; I=4.732397001373103; P=11.97410097664093 }
    13 = @{ B=34.36056989042019; C=26.74867764636883; D=4.535521992404589; E=17.30423646035947; F=44.55714983913662; G=2.037103107177901; H=6.5172900259235; I=4.598145302969217; P=11.65383825406395 }
    14 = @{ B=33.34101321820851; C=25.82294158201565; D=4.309114174977028; E=13.23754183903135; F=42.3561915673402; G=2.040608279471432; H=7.264252431868699; I=4.464898670560597; P=11.49213779797367 }
    15 = @{ B=32.9648862501232; C=25.51230040168827; D=4.257088291712681; E=12.24706028960986; F=41.66581721299364; G=2.042053549585688; H=7.433725654074647; I=4.411213507189735; P=11.46914446365359 }
    16 = @{ B=31.87608673764142; C=24.67704261173039; D=4.270210178006453; E=11.94105676579287; F=40.52409631443208; G=2.047802731150009; H=7.12150906140563; I=4.199051954891138; P=11.6109875496389 }
    17 = @{ B=31.59081828550146; C=24.54023630133161; D=4.392134655494539; E=13.92327043899569; F=40.8602284027189; G=2.050415076146294; H=6.426382675262424; I=4.104955505081555; P=11.7975146447877 }
    18 = @{ B=32.01266221312296; C=24.99574270966277; D=4.659682672957626; E=18.43728141267208; F=42.57297928984077; G=2.050369575357616; H=5.417158518650544; I=4.10792768057442; P=12.08228754755286 }
    19 = @{ B=32.91428265727141; C=25.92709831785335; D=5.046921090650943; E=24.73815109591906; F=45.11864271034842; G=2.04795465780415; H=4.505350713937855; I=4.201341048233777; P=12.44996205798263 }
    20 = @{ B=35.18532413088438; C=28.05782744234055; D=5.75501748214636; E=35.0635044301351; F=50.01271585242392; G=2.038744722034453; H=4.378823149754055; I=4.545497348120654; P=13.05914361158725 }
    21 = @{ B=37.60410681397683; C=30.04120099125489; D=5.965210850010251; E=38.14978486046395; F=53.24885792036237; G=2.026048842488905; H=4.886079521493445; I=5.01815482679849; P=13.03976532144044 }
    22 = @{ B=39.08938018609194; C=31.20984196780294; D=6.081779563242853; E=39.6612641198891; F=55.20015203280625; G=2.01799549018383; H=5.192950727496313; I=5.317053975957242; P=13.01478655940469 }
    23 = @{ B=38.3363028547522; C=30.57376622537857; D=6.034294527975277; E=38.86161477420278; F=54.2361464001029; G=2.022226197195321; H=5.030911046414971; I=5.158508503327936; P=13.04675674622192 }
    24 = @{ B=35.31861004282749; C=28.11630933623611; D=5.827425069724348; E=35.73341028950765; F=50.37637003397209; G=2.038380818390848; H=4.411048834133488; I=4.555896700817713; P=13.14098249754714 }
    25 = @{ B=31.81816931407905; C=25.29905695135953; D=5.596600415070315; E=32.14356109270349; F=46.01060943009153; G=2.055980476337377; H=3.730515494968214; I=3.901720770525965; P=13.2493814597989 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}

Write-Output "done"